$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap B9 and C9 values
$b9 = $ws.Range("B9").Value2
$c9 = $ws.Range("C9").Value2
$ws.Range("B9").Value = $c9
$ws.Range("C9").Value = $b9

# Fix B29 value (Leadership content -> Sustainability)
$ws.Range("B29").Value = "Sustainability"

# Update view: selection moves to B9 (freeze pane top-left cell resets to A2)
$ws.Activate()
$ws.Range("B9").Select()
